# Generate Report for Handoff
#
# The handoff status moves from "Handed back: in sync with en-US" to
# "Ready for handoff" everywhere it is reported (Overview rollup columns
# for each locale, plus the per-locale "Status" column), and the
# corresponding handoff timestamps are refreshed to the moment the new
# report was generated. The "Status"/date columns that used to hold the
# long "Handed back..." text are narrowed now that the new text is much
# shorter.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet: per-locale status + latest handoff xliff generation date ---
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-13 15:13:43"

# --- zh-cn sheet: Status + Latest Handoff Datetime ---
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-13 15:13:35"

# --- de-de sheet: Status + Latest Handoff Datetime ---
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-13 15:13:43"

# Re-fit the columns that used to carry the long status text so they
# shrink down to the width of the new, much shorter "Ready for handoff".
$wsOverview.Columns("E:E").ColumnWidth = 16.3
$wsOverview.Columns("F:F").ColumnWidth = 16.3
$wsZhCn.Columns("C:C").ColumnWidth = 16.3
$wsDeDe.Columns("C:C").ColumnWidth = 16.3
